$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")

$ws.Range("Q2").NumberFormat = "@"
$ws.Range("Q2").Value = "51528308"

$ws.Range("Q3").NumberFormat = "@"
$ws.Range("Q3").Value = "51528312"

$ws.Range("R3").NumberFormat = "@"
$ws.Range("R3").Value = "51528316"

$ws.Range("AD3").NumberFormat = "@"
$ws.Range("AD3").Value = "01-26-2022"

$ws.Range("Q4").NumberFormat = "@"
$ws.Range("Q4").Value = "51528318"
